$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.967.14"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "1.772.90"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "220.67"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "30.91"
$ws.Range("E8").Value = "  -7.03%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "0.0705"
$ws.Range("E10").Value = "  +4.60%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "2.026.45"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "1.774.89"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "10.48"
$ws.Range("E14").Value = "  -5.42%  "
$ws.Range("D15").Value = "0.622"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "33.962.31"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "67.62"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").Value = "242.86"
$ws.Range("E19").Value = "  -5.44%  "
$ws.Range("D20").Value = "0.0₃0772"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "10.52"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "157.73"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "16.33"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "6.98"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "3.69"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "3.49"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").Value = "1.394.93"
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "0.630"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "78.77"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "0.0491"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").Value = "5.85"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "1.923.64"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "103.93"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "11.81"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "0.0₆0119"
$ws.Range("E51").Value = "  -1.23%  "
